# Add season-record columns (Wins, Losses, Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Copy the existing header formatting (bold font, centered, bordered)
# from A1 onto the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows ----------------------------------------------------------
# Every player (rows 2-43) shares the same 2001 Los Angeles Dodgers
# season record: 86 wins, 76 losses, 0 ties.
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 86  # AD
    $ws.Cells.Item($r, 31).Value = 76  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
